# GridCal shortcut script
# Applies three logical changes to the Grid.xlsx workbook:
#   1. config!Name  "Grid" -> "Bloc energy grid"
#   2. bus sheet: shift every bus (x, y) diagram position by (-367, -183)
#   3. branch sheet: "length" becomes "1.0" (was "1"); the G/B admittance
#      value's string form becomes "1e-20" (was "1.0000000000000001e-20")

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, $text)
    # Force the cell to be written as a text/string value (matches the
    # source tool's export, which always writes these columns as text)
    # rather than letting Excel auto-coerce numeric-looking text to a
    # number. Resetting the NumberFormat back to General afterwards keeps
    # the cell's style identical to a plain, unformatted cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# 1. config sheet: rename the grid
$configWs = $wb.Worksheets.Item("config")
$configWs.Range("C4").Value = "Bloc energy grid"

# 2. bus sheet: shift the (x, y) diagram coordinates of every bus
$busWs = $wb.Worksheets.Item("bus")

$dx = -367.0
$dy = -183.0

$busRows = 2..8
foreach ($r in $busRows) {
    $xCell = $busWs.Cells.Item($r, 9)   # column I = x
    $yCell = $busWs.Cells.Item($r, 10)  # column J = y

    $oldX = [double]$xCell.Value2
    $oldY = [double]$yCell.Value2

    $newX = $oldX + $dx
    $newY = $oldY + $dy

    $newXText = $newX.ToString("0.0")
    $newYText = $newY.ToString("0.0")

    Set-TextValue $xCell $newXText
    Set-TextValue $yCell $newYText
}

# 3. branch sheet: normalize the "length" and admittance text values
$branchWs = $wb.Worksheets.Item("branch")

$gbText = "1e-20"
$lenText = "1.0"

$branchRows = 2..10
foreach ($r in $branchRows) {
    $gCell = $branchWs.Cells.Item($r, 11)  # column K = G
    $bCell = $branchWs.Cells.Item($r, 12)  # column L = B
    $lenCell = $branchWs.Cells.Item($r, 13) # column M = length

    Set-TextValue $gCell $gbText
    Set-TextValue $bCell $gbText
    Set-TextValue $lenCell $lenText
}
